# toyexample_TEMPTED.xlsx sign-flip fix
#
# The TEMPTED decomposition has a sign indeterminacy between the subject
# loadings ("A.hat", Component 1 & 2) and the time-loading functions
# ("Phi.hat", Component 1 & 2). The commit flips the sign of those two
# columns in both sheets so the factorization is presented consistently
# (Component 3 / column C of Phi.hat is left untouched).

$wb = $excel.ActiveWorkbook

# --- "A.hat": negate Component 1 (col B) and Component 2 (col C), rows 2-61 ---
$wsA = $wb.Worksheets.Item("A.hat")
$rngA = $wsA.Range("B2:C61")
$valA = $rngA.Value()
$rowsA = $valA.GetLength(0)
$colsA = $valA.GetLength(1)
for ($i = 1; $i -le $rowsA; $i++) {
    for ($j = 1; $j -le $colsA; $j++) {
        $valA[$i, $j] = -1 * $valA[$i, $j]
    }
}
$rngA.Value = $valA

# --- "Phi.hat": negate Component 1 (col A) and Component 2 (col B), rows 2-102 ---
$wsPhi = $wb.Worksheets.Item("Phi.hat")
$rngPhi = $wsPhi.Range("A2:B102")
$valPhi = $rngPhi.Value()
$rowsPhi = $valPhi.GetLength(0)
$colsPhi = $valPhi.GetLength(1)
for ($i = 1; $i -le $rowsPhi; $i++) {
    for ($j = 1; $j -le $colsPhi; $j++) {
        $valPhi[$i, $j] = -1 * $valPhi[$i, $j]
    }
}
$rngPhi.Value = $valPhi
